$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old column G ("img": a single course image path) is removed entirely.
# Everything from former column H onward (images/vdoUrl/desc) shifts left.
$ws.Range("G1").EntireColumn.Delete()

# Re-point the two hyperlink cells that used to live in column I (now column H)
$ws.Range("H3").Hyperlinks.Delete()
$ws.Range("H5").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("H3"), "http://www.youtube.com/")
$ws.Hyperlinks.Add($ws.Range("H5"), "http://www.youtube.com/")

# Column widths (characters)
$ws.Columns.Item(1).ColumnWidth = 9.140625
$ws.Columns.Item(2).ColumnWidth = 27.5703125
$ws.Columns.Item(3).ColumnWidth = 17.42578125
$ws.Columns.Item(4).ColumnWidth = 9.140625
$ws.Columns.Item(5).ColumnWidth = 9.140625
$ws.Columns.Item(6).ColumnWidth = 9.140625
$ws.Columns.Item(7).ColumnWidth = 47.42578125
$ws.Columns.Item(8).ColumnWidth = 20.28515625
$ws.Columns.Item(9).ColumnWidth = 63.28515625

# Wrap text across the used grid (header + data), cell-by-cell run so we
# don't manufacture style-only entries for the (legitimately) empty cells
$ws.Range("A1:I1").WrapText = $true
$ws.Range("A2:G2").WrapText = $true
$ws.Range("I2").WrapText = $true
$ws.Range("A3:I3").WrapText = $true
$ws.Range("A4:G4").WrapText = $true
$ws.Range("I4").WrapText = $true
$ws.Range("A5:I5").WrapText = $true

# Row heights to fit the wrapped, multi-line descriptions
$ws.Rows.Item(2).RowHeight = 45
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 30

# Make sure the header font is still explicitly Calibri (keeps font metadata in sync)
$ws.Range("A1:I1").Font.Name = "Calibri"

# View: zoom to 70% and move the selection like the saved workbook shows
$excel.ActiveWindow.Zoom = 70
$ws.Range("H10").Select()
